$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Val)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Val
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.968.89'
Set-TextValue $ws.Range('E2') '  -0.19%  '
Set-TextValue $ws.Range('D3') '1.868.54'
Set-TextValue $ws.Range('E3') '  -2.66%  '
Set-TextValue $ws.Range('E4') '  +0.06%  '
Set-TextValue $ws.Range('D5') '319.28'
Set-TextValue $ws.Range('E5') '  -2.64%  '
Set-TextValue $ws.Range('E6') '  +0.04%  '
Set-TextValue $ws.Range('D7') '0.5082'
Set-TextValue $ws.Range('E7') '  -3.16%  '
Set-TextValue $ws.Range('D8') '0.3933'
Set-TextValue $ws.Range('E8') '  -2.92%  '
Set-TextValue $ws.Range('D9') '0.08177'
Set-TextValue $ws.Range('E9') '  -3.44%  '
Set-TextValue $ws.Range('D10') '41.99'
Set-TextValue $ws.Range('E10') '  -2.13%  '
Set-TextValue $ws.Range('D11') '1.090'
Set-TextValue $ws.Range('E11') '  -3.20%  '
Set-TextValue $ws.Range('D12') '22.80'
Set-TextValue $ws.Range('E12') '  +2.59%  '
Set-TextValue $ws.Range('D13') '1.859.45'
Set-TextValue $ws.Range('E13') '  -3.12%  '
Set-TextValue $ws.Range('D14') '6.270'
Set-TextValue $ws.Range('E14') '  -1.34%  '
Set-TextValue $ws.Range('D15') '7.166'
Set-TextValue $ws.Range('E15') '  -2.79%  '
Set-TextValue $ws.Range('D17') '91.74'
Set-TextValue $ws.Range('E17') '  -4.71%  '
Set-TextValue $ws.Range('E18') '  -2.68%  '
Set-TextValue $ws.Range('D19') '0.06398'
Set-TextValue $ws.Range('E19') '  -4.97%  '
Set-TextValue $ws.Range('D20') '17.85'
Set-TextValue $ws.Range('E20') '  -2.01%  '
Set-TextValue $ws.Range('D21') '1.000'
Set-TextValue $ws.Range('E21') '  -0.05%  '
Set-TextValue $ws.Range('D22') '29.967.42'
Set-TextValue $ws.Range('E22') '  -0.13%  '
Set-TextValue $ws.Range('D23') '5.809'
Set-TextValue $ws.Range('E23') '  -4.00%  '
Set-TextValue $ws.Range('D24') '11.08'
Set-TextValue $ws.Range('E24') '  -1.52%  '
Set-TextValue $ws.Range('D25') '2.164'
Set-TextValue $ws.Range('E25') '  -3.02%  '
Set-TextValue $ws.Range('D26') '2.078.07'
Set-TextValue $ws.Range('E26') '  -2.92%  '
Set-TextValue $ws.Range('D27') '160.97'
Set-TextValue $ws.Range('E27') '  +0.59%  '
Set-TextValue $ws.Range('D28') '20.91'
Set-TextValue $ws.Range('D29') '2.217'
Set-TextValue $ws.Range('E29') '  -9.64%  '
Set-TextValue $ws.Range('D30') '127.17'
Set-TextValue $ws.Range('E30') '  -1.79%  '
Set-TextValue $ws.Range('D31') '1.064'
Set-TextValue $ws.Range('E31') '  -1.57%  '
Set-TextValue $ws.Range('D32') '0.1034'
Set-TextValue $ws.Range('E32') '  -2.54%  '
Set-TextValue $ws.Range('E33') '  -3.32%  '
Set-TextValue $ws.Range('E34') '  +1.84%  '
Set-TextValue $ws.Range('D35') '0.02423'
Set-TextValue $ws.Range('E35') '  -3.87%  '
Set-TextValue $ws.Range('D36') '5.194'
Set-TextValue $ws.Range('E36') '  -0.06%  '
Set-TextValue $ws.Range('D37') '0.06323'
Set-TextValue $ws.Range('E37') '  -3.98%  '
Set-TextValue $ws.Range('D38') '0.2135'
Set-TextValue $ws.Range('E38') '  -4.14%  '
Set-TextValue $ws.Range('D39') '1.168'
Set-TextValue $ws.Range('E39') '  -5.56%  '
Set-TextValue $ws.Range('D40') '8.459'
Set-TextValue $ws.Range('E40') '  -5.90%  '
Set-TextValue $ws.Range('D41') '0.6285'
Set-TextValue $ws.Range('E41') '  -4.08%  '
Set-TextValue $ws.Range('B42') 'TrustWalletToken'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D42') '1.203'
Set-TextValue $ws.Range('E42') '  -3.67%  '
Set-TextValue $ws.Range('B43') 'Aptos'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D43') '11.21'
Set-TextValue $ws.Range('E43') '  -3.65%  '
Set-TextValue $ws.Range('E44') '  +0.00%  '
Set-TextValue $ws.Range('B45') 'Decentraland'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D45') '0.5872'
Set-TextValue $ws.Range('E45') '  -5.05%  '
Set-TextValue $ws.Range('B46') 'EnergySwap'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '12.87'
Set-TextValue $ws.Range('E46') '  -2.32%  '
Set-TextValue $ws.Range('D47') '3.635'
Set-TextValue $ws.Range('E47') '  -3.28%  '
Set-TextValue $ws.Range('D48') '1.991'
Set-TextValue $ws.Range('E48') '  -3.75%  '
Set-TextValue $ws.Range('D49') '122.13'
Set-TextValue $ws.Range('E49') '  -2.88%  '
Set-TextValue $ws.Range('D50') '1.198'
Set-TextValue $ws.Range('E50') '  -3.54%  '
Set-TextValue $ws.Range('D51') '1.118'
